# Fix project gantt chart
#
# The "EstimatedEffortHours" (col C) and "Progress" (col E) columns are
# removed from the projects table, leaving ID | ProjectName | Dependencies.
# The remaining "Dependencies" column (was col D) slides left into col C,
# taking its column width/formatting with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete "Progress" (E) first so the "EstimatedEffortHours" (C) delete
# below doesn't shift its position.
$ws.Range("E1:E3").EntireColumn.Delete()

# Delete "EstimatedEffortHours" (C); "Dependencies" (old D) slides into C.
$ws.Range("C1:C3").EntireColumn.Delete()

# Make sure the header text is exactly right after the shift.
$ws.Range("C1").Value = "Dependencies"

# Shrink the table to the new A1:C3 extent (drops the two removed
# columns from the table definition and refreshes the autofilter range).
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C3"))

# Match the saved selection from the source edit.
[void]$ws.Range("C5").Select()
